$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I (I0) and J (IF), matching style of existing header row (B1:H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill I2:I31 with 1, and J2:J31 mirrors H2:H31
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
}
